$d = $word.ActiveDocument

$pairs = @(
  @("444÷3=148, 0", "295÷9=32, 7"),
  @("621÷8=77, 5", "354÷9=39, 3"),
  @("842÷8=105, 2", "448÷7=64, 0"),
  @("917÷6=152, 5", "765÷7=109, 2"),
  @("115÷9=12, 7", "998÷5=199, 3"),
  @("556÷3=185, 1", "811÷7=115, 6"),
  @("439÷8=54, 7", "950÷6=158, 2"),
  @("906÷9=100, 6", "825÷6=137, 3"),
  @("338÷9=37, 5", "921÷5=184, 1"),
  @("797÷4=199, 1", "925÷7=132, 1"),
  @("622÷4=155, 2", "753÷5=150, 3"),
  @("665÷3=221, 2", "739÷8=92, 3"),
  @("744÷5=148, 4", "337÷2=168, 1"),
  @("520÷8=65, 0", "654÷5=130, 4"),
  @("460÷7=65, 5", "236÷5=47, 1"),
  @("382÷2=191, 0", "132÷7=18, 6"),
  @("193÷6=32, 1", "178÷5=35, 3"),
  @("700÷9=77, 7", "473÷3=157, 2"),
  @("398÷6=66, 2", "400÷7=57, 1"),
  @("542÷8=67, 6", "775÷4=193, 3"),
  @("190÷3=63, 1", "578÷3=192, 2"),
  @("627÷4=156, 3", "323÷5=64, 3"),
  @("667÷3=222, 1", "327÷5=65, 2"),
  @("375÷9=41, 6", "533÷7=76, 1"),
  @("528÷5=105, 3", "356÷4=89, 0")
)

foreach ($pair in $pairs) {
  $old = $pair[0]
  $new = $pair[1]
  $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
